# Auto update Excel log: append 6 new PRESENCE_DETECTED mmWave rows
# (2026-01-31, 21:34:59 .. 21:35:49) to the "mmWave" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

$rows = @(
    @("2026-01-31", "21:34:59", "21:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-31", "21:35:07", "21:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-31", "21:35:17", "21:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-31", "21:35:28", "21:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-31", "21:35:38", "21:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-31", "21:35:49", "21:00", "Living Room", "PRESENCE_DETECTED", "Active")
)

$startRow = 8

# Columns A-C hold date/time-looking text (e.g. "2026-01-31", "21:00") that
# must stay plain text, matching every existing row on this sheet - format
# those cells as Text first so Excel doesn't coerce them into date/time
# serial numbers.
$endRow = $startRow + $rows.Count - 1
$dateTimeRange = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 3))
$dateTimeRange.NumberFormat = "@"

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}
